$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*, *") {
        $parts = $val -split ", "
        $reversed = $parts[($parts.Length - 1)..0]
        $cell.Value2 = $reversed -join ", "
    }
}
